$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numbers formatted as plain text in the source data (e.g.
# "26.073.73", "1.003", "0.00001076"). Setting .Value directly would let Excel
# auto-coerce numeric-looking strings into real floating point numbers and lose
# the exact text representation, so we prefix with a literal apostrophe to force
# text entry, matching how the source workbook stores these cells.

$ws.Range("D2").Value = '''26.073.73'
$ws.Range("E2").Value = '  +5.72%  '

$ws.Range("D3").Value = '''1.719.61'
$ws.Range("E3").Value = '  +3.79%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''332.12'
$ws.Range("E5").Value = '  +3.69%  '

$ws.Range("D6").Value = '''0.9986'
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '''0.3693'
$ws.Range("E7").Value = '  +1.55%  '

$ws.Range("D8").Value = '''49.51'
$ws.Range("E8").Value = '  +5.66%  '

$ws.Range("D9").Value = '''0.3341'
$ws.Range("E9").Value = '  +2.54%  '

$ws.Range("D10").Value = '''1.187'
$ws.Range("E10").Value = '  +5.00%  '

$ws.Range("D11").Value = '''0.07479'
$ws.Range("E11").Value = '  +6.21%  '

$ws.Range("D12").Value = '''0.9988'
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("D13").Value = '''6.280'
$ws.Range("E13").Value = '  +5.06%  '

$ws.Range("D14").Value = '''20.09'
$ws.Range("E14").Value = '  +3.08%  '

$ws.Range("D15").Value = '''6.918'
$ws.Range("E15").Value = '  +4.58%  '

$ws.Range("D16").Value = '''1.715.27'
$ws.Range("E16").Value = '  +3.21%  '

$ws.Range("D17").Value = '''0.00001076'
$ws.Range("E17").Value = '  +3.00%  '

$ws.Range("D18").Value = '''0.06639'
$ws.Range("E18").Value = '  +0.45%  '

$ws.Range("D19").Value = '''81.99'
$ws.Range("E19").Value = '  +4.16%  '

$ws.Range("D20").Value = '''0.9989'
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").Value = '''16.40'
$ws.Range("E21").Value = '  +4.36%  '

$ws.Range("D22").Value = '''6.081'
$ws.Range("E22").Value = '  +2.62%  '

$ws.Range("D23").Value = '''13.00'
$ws.Range("E23").Value = '  +2.98%  '

$ws.Range("D24").Value = '''26.031.90'
$ws.Range("E24").Value = '  +5.62%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '''2.469'
$ws.Range("E26").Value = '  +3.34%  '

$ws.Range("D27").Value = '''150.25'
$ws.Range("E27").Value = '  +1.69%  '

$ws.Range("D28").Value = '''19.27'
$ws.Range("E28").Value = '  +3.68%  '

$ws.Range("E29").Value = '  +8.70%  '

$ws.Range("D30").Value = '''1.904.82'
$ws.Range("E30").Value = '  +3.21%  '

$ws.Range("D31").Value = '''129.39'
$ws.Range("E31").Value = '  +3.26%  '

$ws.Range("D32").Value = '''4.090'
$ws.Range("E32").Value = '  +0.35%  '

$ws.Range("D33").Value = '''5.949'

$ws.Range("E34").Value = '  +0.81%  '

$ws.Range("D35").Value = '''1.715'
$ws.Range("E35").Value = '  +1.63%  '

$ws.Range("D36").Value = '''12.90'
$ws.Range("E36").Value = '  +4.77%  '

$ws.Range("D37").Value = '''5.347'
$ws.Range("E37").Value = '  +2.81%  '

$ws.Range("D38").Value = '''1.264'
$ws.Range("E38").Value = '  -0.68%  '

$ws.Range("D39").Value = '''0.06199'
$ws.Range("E39").Value = '  +2.75%  '

$ws.Range("D40").Value = '''0.02289'
$ws.Range("E40").Value = '  +2.54%  '

$ws.Range("D42").Value = '''8.529'
$ws.Range("E42").Value = '  +4.24%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.6168'
$ws.Range("E43").Value = '  +4.22%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''14.38'
$ws.Range("E44").Value = '  +13.28%  '

$ws.Range("D45").Value = '''0.9990'
$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").Value = '''3.831'
$ws.Range("E46").Value = '  -0.43%  '

$ws.Range("D47").Value = '''0.5885'
$ws.Range("E47").Value = '  +4.92%  '

$ws.Range("E48").Value = '  +2.88%  '

$ws.Range("D49").Value = '''2.021'
$ws.Range("E49").Value = '  +3.48%  '

$ws.Range("D50").Value = '''0.07263'
$ws.Range("E50").Value = '  +4.21%  '

$ws.Range("D51").Value = '''77.11'
$ws.Range("E51").Value = '  +3.50%  '
